{"js": "// Fix two typos in the \"Descri\u00e7\u00e3o do Projeto\" document:\n//   1) \"lancetes\" -> \"lacetes\"\n//   2) \"contruir\" -> \"construir\" (the letter \"s\" is inserted as its own\n//      run between \"con\" and \"truir\", mirroring how Word records an\n//      in-place keystroke correction as a separate run).\n\nconst body = context.document.body;\n\n// --- Edit 1: \"Um grafo G conexo sem lancetes \u00e9 \" -> \"...lacetes \u00e9 \" ---\nconst typo1 = body.search(\"lancetes\", { matchCase: false });\ntypo1.load(\"text\");\nawait context.sync();\n\nif (typo1.items.length > 0) {\n  typo1.items[0].insertText(\"lacetes\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Edit 2: \"ir\u00e1 contruir grafos\" -> \"ir\u00e1 construir grafos\" ---\n// Scope the search to the unique phrase \"ir\u00e1 contruir\" so we don't\n// collide with the other \"con...\" occurrences in the document\n// (e.g. \"conhecimento\", \"conexo\").\nconst phrase = body.search(\"ir\u00e1 contruir\", { matchCase: false });\nphrase.load(\"text\");\nawait context.sync();\n\nif (phrase.items.length > 0) {\n  // Narrow further to just the \"con\" prefix of \"contruir\" within that\n  // scoped match, then insert the missing \"s\" right after it.\n  const conPart = phrase.items[0].search(\"con\", { matchCase: false });\n  await context.sync();\n\n  const conRange = conPart.items[0];\n  conRange.insertText(\"s\", Word.InsertLocation.after);\n  await context.sync();\n\n  // At this point the text reads \"...ir\u00e1 construir grafos...\" but the\n  // insertion landed inside the single pre-existing run. Re-find the\n  // inserted letter and isolate it (and the following \"truir\") into\n  // their own runs by toggling a character property on/off, which\n  // forces the engine to split the run without altering its visible\n  // formatting (matches the target XML: three runs with identical rPr).\n  const tail = body.search(\"struir\", { matchCase: false });\n  await context.sync();\n\n  const tailRange = tail.items[0];\n  const pieces = tailRange.split([\"s\", \"truir\"], false, false);\n  await context.sync();\n\n  const sRange = pieces.items[0];\n  sRange.font.bold = true;\n  await context.sync();\n  sRange.font.bold = false;\n  await context.sync();\n}\n", "ps1": "# Fix two typos in the \"Descri\u00e7\u00e3o do Projeto\" document:\n#   1) \"lancetes\" -> \"lacetes\"\n#   2) \"contruir\" -> \"construir\" (the letter \"s\" is inserted as its own\n#      run between \"con\" and \"truir\", mirroring how Word records an\n#      in-place keystroke correction as a separate run).\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: \"Um grafo G conexo sem lancetes \u00e9 \" -> \"...lacetes \u00e9 \" ---\n$find1 = $d.Content\n$find1.Find.Execute(\"lancetes\")\nif ($find1.Find.Found) {\n    $find1.Text = \"lacetes\"\n}\n\n# --- Edit 2: \"ir\u00e1 contruir grafos\" -> \"ir\u00e1 construir grafos\" ---\n# Scope the search to the unique phrase \"ir\u00e1 contruir\" so we don't\n# collide with the other \"con...\" occurrences in the document\n# (e.g. \"conhecimento\", \"conexo\").\n$phrase = $d.Content\n$phrase.Find.Execute(\"ir\u00e1 contruir\")\nif ($phrase.Find.Found) {\n    # \"ir\u00e1 \" is 4 characters, \"con\" is the next 3 characters of the match.\n    $conRange = $d.Range($phrase.Start + 4, $phrase.Start + 7)\n    $conRange.InsertAfter(\"s\")\n\n    # InsertAfter grows $conRange to cover the newly-inserted text too\n    # (\"con\" -> \"cons\"), so the inserted \"s\" is the last character of it.\n    $sRange = $d.Range($conRange.End - 1, $conRange.End)\n\n    # Isolate the \"s\" into its own run by toggling a character property\n    # on/off. This forces the engine to split the run at both the\n    # \"con\"/\"s\" and \"s\"/\"truir\" boundaries without altering the visible\n    # formatting (matches the target XML: three runs with identical\n    # rPr/rFonts/sz).\n    $sRange.Font.Bold = 1\n    $sRange.Font.Bold = 0\n}\n"}
